$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's data was stored as an Excel Table ("Table1" over A1:I19).
# Convert it back to a normal range (removes xl/tables/table1.xml,
# the tableParts reference, and the table-driven dxfs/header border).
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# Update the Warranty Start / Warranty End dates for the 3rd data row
# (row 4) with new test values - these become two new shared strings.
# Force them in as text (matching the existing text-stored dates in the
# rest of the column) instead of letting Excel auto-convert them to
# date serial numbers, then restore the Normal style so no stray
# number-format style lingers on the cell.
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "06-09-2026"
$ws.Range("F4").Style = "Normal"

$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "09-07-2029"
$ws.Range("G4").Style = "Normal"

# Matches the updated cursor position recorded in the saved workbook.
$ws.Range("G12").Select()
